$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1664
$ws.Range("F3").Value = 868
$ws.Range("G3").Value = 100
$ws.Range("F4").Value = 275
$ws.Range("F5").Value = 84
$ws.Range("F6").Value = 1189
$ws.Range("F7").Value = 813
$ws.Range("F8").Value = 836
$ws.Range("F9").Value = 1542
$ws.Range("F10").Value = 310
$ws.Range("F11").Value = 1063
$ws.Range("F14").Value = 206
$ws.Range("F16").Value = 517
$ws.Range("F17").Value = 70
$ws.Range("F18").Value = 42
$ws.Range("F22").Value = 582
$ws.Range("F23").Value = 587
$ws.Range("F24").Value = 59
$ws.Range("F25").Value = 10
$ws.Range("F26").Value = 784
$ws.Range("F27").Value = 263
$ws.Range("F28").Value = 198
$ws.Range("F30").Value = 377

# ---- Sheet "演出" (performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1040
$ws.Range("F5").Value = 284

# ---- Sheet "全部类型" (all types, combined view) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1664
$ws.Range("F5").Value = 868
$ws.Range("G5").Value = 100
$ws.Range("F6").Value = 275
$ws.Range("F7").Value = 1040
$ws.Range("F8").Value = 84
$ws.Range("F9").Value = 1189
$ws.Range("F10").Value = 813
$ws.Range("F11").Value = 836
$ws.Range("F12").Value = 1542
$ws.Range("F13").Value = 310
$ws.Range("F14").Value = 1063
$ws.Range("F17").Value = 206
$ws.Range("F19").Value = 517
$ws.Range("F20").Value = 70
$ws.Range("F21").Value = 42
$ws.Range("F24").Value = 284
$ws.Range("F30").Value = 582
$ws.Range("F31").Value = 587
$ws.Range("F32").Value = 59
$ws.Range("F33").Value = 10
$ws.Range("F34").Value = 784
$ws.Range("F35").Value = 263
$ws.Range("F37").Value = 198
$ws.Range("F43").Value = 377
